# [Kadastro App] Yeni kayit eklendi: 2967
# Adds a new record row (row 35) to the "Kayitlar" (master records) sheet
# and to the matching "Erdemli" district sheet, mirroring the values that
# already exist for that district/job-type in the workbook.

$wb = $excel.ActiveWorkbook

$newRow = @{
    A = "2967"
    B = "2025-09-10"
    C = "Erdemli"
    D = "1"
    E = "ÇAP"
    F = "AYHAN KARADAYI (K.Teknisyeni)"
}

$targetSheets = @("Kayitlar", "Erdemli")

foreach ($sheetName in $targetSheets) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Row directly after the current last row of data (row 34 -> 35).
    $rowIndex = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row + 1

    $rowRange = $ws.Range("A" + $rowIndex + ":F" + $rowIndex)

    # Force the new cells to be stored as text (matching the rest of the
    # sheet, which keeps numeric-looking values like the record id and
    # parcel count as text) instead of being auto-converted to numbers
    # or dates.
    $rowRange.NumberFormat = "@"

    $ws.Range("A" + $rowIndex).Value = $newRow.A
    $ws.Range("B" + $rowIndex).Value = $newRow.B
    $ws.Range("C" + $rowIndex).Value = $newRow.C
    $ws.Range("D" + $rowIndex).Value = $newRow.D
    $ws.Range("E" + $rowIndex).Value = $newRow.E
    $ws.Range("F" + $rowIndex).Value = $newRow.F

    # Drop the explicit text-number-format style again so the new cells
    # stay styleless (same as the rest of the sheet) while remaining
    # text-typed values.
    $rowRange.ClearFormats()
}
